$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing A3 timestamp value (tiny float precision fix)
$ws.Range("A3").Value = 45804.45374693287

# Add new row 4 with the updated price entry
$ws.Range("A4").Value = 45805.39378755033
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B4").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C4").Value = "1Kg"
$ws.Range("D4").Value = "15,41€"
